# Update "want to go" (想去人数) counts in column F on the "展览" sheet
# and the "全部类型" sheet, mirroring the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F5").Value  = 1837
$ws1.Range("F6").Value  = 472
$ws1.Range("F7").Value  = 167
$ws1.Range("F8").Value  = 166
$ws1.Range("F9").Value  = 2428
$ws1.Range("F10").Value = 143
$ws1.Range("F11").Value = 78
$ws1.Range("F12").Value = 163
$ws1.Range("F13").Value = 1472
$ws1.Range("F15").Value = 37
$ws1.Range("F19").Value = 180
$ws1.Range("F25").Value = 40
$ws1.Range("F26").Value = 1517
$ws1.Range("F28").Value = 379
$ws1.Range("F29").Value = 341
$ws1.Range("F30").Value = 187
$ws1.Range("F31").Value = 289
$ws1.Range("F32").Value = 379

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value  = 1837
$ws4.Range("F7").Value  = 472
$ws4.Range("F8").Value  = 167
$ws4.Range("F9").Value  = 166
$ws4.Range("F10").Value = 2428
$ws4.Range("F11").Value = 143
$ws4.Range("F12").Value = 78
$ws4.Range("F13").Value = 163
$ws4.Range("F14").Value = 1472
$ws4.Range("F16").Value = 37
$ws4.Range("F20").Value = 180
$ws4.Range("F26").Value = 40
$ws4.Range("F27").Value = 1517
$ws4.Range("F29").Value = 379
$ws4.Range("F30").Value = 341
$ws4.Range("F31").Value = 187
$ws4.Range("F32").Value = 289
$ws4.Range("F33").Value = 379
